# Recompute NATMI ligand-receptor output with updated TPM-based values.
# The "Inflammatory-Mac" target-cluster rows are dropped, the remaining
# FAPs/MuSCs x {ECs,FAPs,MuSCs,Neutrophils,Resolving-Mac} combinations
# (rows 2-11) get refreshed numeric results, and the now-unused
# "Inflammatory-Mac" shared string disappears automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'FAPs'
$ws.Range('B2').Value = 'Lgi1'
$ws.Range('C2').Value = 'Adam23'
$ws.Range('D2').Value = 'ECs'
$ws.Range('E2').Value = 1.0
$ws.Range('F2').Value = 0.3333333333333333
$ws.Range('G2').Value = 0.002491666666666667
$ws.Range('H2').Value = 0.007475
$ws.Range('I2').Value = 0.1635273785303319
$ws.Range('J2').Value = 0.1635273785303318
$ws.Range('K2').Value = 3.0
$ws.Range('L2').Value = 1.0
$ws.Range('M2').Value = 0.2328916666666666
$ws.Range('N2').Value = 0.6986749999999999
$ws.Range('O2').Value = 0.01421300418632399
$ws.Range('P2').Value = 0.01421300418632399
$ws.Range('Q2').Value = 0.0005802884027777777
$ws.Range('R2').Value = 0.005222595625
$ws.Range('S2').Value = 0.002324215315630195
$ws.Range('T2').Value = 0.002324215315630195

$ws.Range('A3').Value = 'FAPs'
$ws.Range('B3').Value = 'Lgi1'
$ws.Range('C3').Value = 'Adam23'
$ws.Range('D3').Value = 'FAPs'
$ws.Range('E3').Value = 1.0
$ws.Range('F3').Value = 0.3333333333333333
$ws.Range('G3').Value = 0.002491666666666667
$ws.Range('H3').Value = 0.007475
$ws.Range('I3').Value = 0.1635273785303319
$ws.Range('J3').Value = 0.1635273785303318
$ws.Range('K3').Value = 3.0
$ws.Range('L3').Value = 1.0
$ws.Range('M3').Value = 14.57672866666667
$ws.Range('N3').Value = 43.730186
$ws.Range('O3').Value = 0.8895943273864486
$ws.Range('P3').Value = 0.8895943273864487
$ws.Range('Q3').Value = 0.03632034892777778
$ws.Range('R3').Value = 0.32688314035
$ws.Range('S3').Value = 0.1454730283129597
$ws.Range('T3').Value = 0.1454730283129597

$ws.Range('A4').Value = 'FAPs'
$ws.Range('B4').Value = 'Lgi1'
$ws.Range('C4').Value = 'Adam23'
$ws.Range('D4').Value = 'MuSCs'
$ws.Range('E4').Value = 1.0
$ws.Range('F4').Value = 0.3333333333333333
$ws.Range('G4').Value = 0.002491666666666667
$ws.Range('H4').Value = 0.007475
$ws.Range('I4').Value = 0.1635273785303319
$ws.Range('J4').Value = 0.1635273785303318
$ws.Range('K4').Value = 3.0
$ws.Range('L4').Value = 1.0
$ws.Range('M4').Value = 1.423796666666667
$ws.Range('N4').Value = 4.27139
$ws.Range('O4').Value = 0.08689202268783405
$ws.Range('P4').Value = 0.08689202268783405
$ws.Range('Q4').Value = 0.003547626694444445
$ws.Range('R4').Value = 0.03192864025
$ws.Range('S4').Value = 0.01420922468533962
$ws.Range('T4').Value = 0.01420922468533962

$ws.Range('A5').Value = 'FAPs'
$ws.Range('B5').Value = 'Lgi1'
$ws.Range('C5').Value = 'Adam23'
$ws.Range('D5').Value = 'Neutrophils'
$ws.Range('E5').Value = 1.0
$ws.Range('F5').Value = 0.3333333333333333
$ws.Range('G5').Value = 0.002491666666666667
$ws.Range('H5').Value = 0.007475
$ws.Range('I5').Value = 0.1635273785303319
$ws.Range('J5').Value = 0.1635273785303318
$ws.Range('K5').Value = 3.0
$ws.Range('L5').Value = 1.0
$ws.Range('M5').Value = 0.136948
$ws.Range('N5').Value = 0.410844
$ws.Range('O5').Value = 0.008357716380185487
$ws.Range('P5').Value = 0.008357716380185487
$ws.Range('Q5').Value = 0.0003412287666666666
$ws.Range('R5').Value = 0.0030710589
$ws.Range('S5').Value = 0.001366715450151747
$ws.Range('T5').Value = 0.001366715450151747

$ws.Range('A6').Value = 'FAPs'
$ws.Range('B6').Value = 'Lgi1'
$ws.Range('C6').Value = 'Adam23'
$ws.Range('D6').Value = 'Resolving-Mac'
$ws.Range('E6').Value = 1.0
$ws.Range('F6').Value = 0.3333333333333333
$ws.Range('G6').Value = 0.002491666666666667
$ws.Range('H6').Value = 0.007475
$ws.Range('I6').Value = 0.1635273785303319
$ws.Range('J6').Value = 0.1635273785303318
$ws.Range('K6').Value = 2.0
$ws.Range('L6').Value = 0.6666666666666666
$ws.Range('M6').Value = 0.01545066666666667
$ws.Range('N6').Value = 0.046352
$ws.Range('O6').Value = 0.0009429293592077716
$ws.Range('P6').Value = 0.0009429293592077717
$ws.Range('Q6').Value = 0.00003849791111111112
$ws.Range('R6').Value = 0.0003464812
$ws.Range('S6').Value = 0.0001541947662505325
$ws.Range('T6').Value = 0.0001541947662505325

$ws.Range('A7').Value = 'MuSCs'
$ws.Range('B7').Value = 'Lgi1'
$ws.Range('C7').Value = 'Adam23'
$ws.Range('D7').Value = 'ECs'
$ws.Range('E7').Value = 2.0
$ws.Range('F7').Value = 0.6666666666666666
$ws.Range('G7').Value = 0.01274533333333334
$ws.Range('H7').Value = 0.03823600000000001
$ws.Range('I7').Value = 0.8364726214696682
$ws.Range('J7').Value = 0.836472621469668
$ws.Range('K7').Value = 3.0
$ws.Range('L7').Value = 1.0
$ws.Range('M7').Value = 0.2328916666666666
$ws.Range('N7').Value = 0.6986749999999999
$ws.Range('O7').Value = 0.01421300418632399
$ws.Range('P7').Value = 0.01421300418632399
$ws.Range('Q7').Value = 0.002968281922222222
$ws.Range('R7').Value = 0.0267145373
$ws.Range('S7').Value = 0.0118887888706938
$ws.Range('T7').Value = 0.0118887888706938

$ws.Range('A8').Value = 'MuSCs'
$ws.Range('B8').Value = 'Lgi1'
$ws.Range('C8').Value = 'Adam23'
$ws.Range('D8').Value = 'FAPs'
$ws.Range('E8').Value = 2.0
$ws.Range('F8').Value = 0.6666666666666666
$ws.Range('G8').Value = 0.01274533333333334
$ws.Range('H8').Value = 0.03823600000000001
$ws.Range('I8').Value = 0.8364726214696682
$ws.Range('J8').Value = 0.836472621469668
$ws.Range('K8').Value = 3.0
$ws.Range('L8').Value = 1.0
$ws.Range('M8').Value = 14.57672866666667
$ws.Range('N8').Value = 43.730186
$ws.Range('O8').Value = 0.8895943273864486
$ws.Range('P8').Value = 0.8895943273864487
$ws.Range('Q8').Value = 0.1857852657662223
$ws.Range('R8').Value = 1.672067391896
$ws.Range('S8').Value = 0.7441212990734889
$ws.Range('T8').Value = 0.7441212990734889

$ws.Range('A9').Value = 'MuSCs'
$ws.Range('B9').Value = 'Lgi1'
$ws.Range('C9').Value = 'Adam23'
$ws.Range('D9').Value = 'MuSCs'
$ws.Range('E9').Value = 2.0
$ws.Range('F9').Value = 0.6666666666666666
$ws.Range('G9').Value = 0.01274533333333334
$ws.Range('H9').Value = 0.03823600000000001
$ws.Range('I9').Value = 0.8364726214696682
$ws.Range('J9').Value = 0.836472621469668
$ws.Range('K9').Value = 3.0
$ws.Range('L9').Value = 1.0
$ws.Range('M9').Value = 1.423796666666667
$ws.Range('N9').Value = 4.27139
$ws.Range('O9').Value = 0.08689202268783405
$ws.Range('P9').Value = 0.08689202268783405
$ws.Range('Q9').Value = 0.01814676311555556
$ws.Range('R9').Value = 0.16332086804
$ws.Range('S9').Value = 0.07268279800249443
$ws.Range('T9').Value = 0.07268279800249443

$ws.Range('A10').Value = 'MuSCs'
$ws.Range('B10').Value = 'Lgi1'
$ws.Range('C10').Value = 'Adam23'
$ws.Range('D10').Value = 'Neutrophils'
$ws.Range('E10').Value = 2.0
$ws.Range('F10').Value = 0.6666666666666666
$ws.Range('G10').Value = 0.01274533333333334
$ws.Range('H10').Value = 0.03823600000000001
$ws.Range('I10').Value = 0.8364726214696682
$ws.Range('J10').Value = 0.836472621469668
$ws.Range('K10').Value = 3.0
$ws.Range('L10').Value = 1.0
$ws.Range('M10').Value = 0.136948
$ws.Range('N10').Value = 0.410844
$ws.Range('O10').Value = 0.008357716380185487
$ws.Range('P10').Value = 0.008357716380185487
$ws.Range('Q10').Value = 0.001745447909333333
$ws.Range('R10').Value = 0.015709031184
$ws.Range('S10').Value = 0.00699100093003374
$ws.Range('T10').Value = 0.006991000930033739

$ws.Range('A11').Value = 'MuSCs'
$ws.Range('B11').Value = 'Lgi1'
$ws.Range('C11').Value = 'Adam23'
$ws.Range('D11').Value = 'Resolving-Mac'
$ws.Range('E11').Value = 2.0
$ws.Range('F11').Value = 0.6666666666666666
$ws.Range('G11').Value = 0.01274533333333334
$ws.Range('H11').Value = 0.03823600000000001
$ws.Range('I11').Value = 0.8364726214696682
$ws.Range('J11').Value = 0.836472621469668
$ws.Range('K11').Value = 2.0
$ws.Range('L11').Value = 0.6666666666666666
$ws.Range('M11').Value = 0.01545066666666667
$ws.Range('N11').Value = 0.046352
$ws.Range('O11').Value = 0.0009429293592077716
$ws.Range('P11').Value = 0.0009429293592077717
$ws.Range('Q11').Value = 0.0001969238968888889
$ws.Range('R11').Value = 0.001772315072
$ws.Range('S11').Value = 0.0007887345929572391
$ws.Range('T11').Value = 0.0007887345929572391

$ws.Range("A12:T13").ClearContents()
